$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Id" column header
$ws.Range("E1").Value = "Id"

# Id values for the existing 3 rows
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3

# Two new student rows
$ws.Range("A5").Value = "Yuleisi"
$ws.Range("B5").Value = "Feliz"
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 20184432
$ws.Range("E5").Value = 4

$ws.Range("A6").Value = "Mary"
$ws.Range("B6").Value = "Villaman"
$ws.Range("C6").Value = 52
$ws.Range("D6").Value = 20184321
$ws.Range("E6").Value = 5

# Update the view: zoom + selection on F1 (past the new data)
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("F1").Select()
